$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column V ("Si autre situation pro") - this
# pushes every existing column from V onward one position to the right
# (V->W, W->X, ... AY->AZ), matching Excel's native "insert column" shift.
$ws.Columns("V:V").Insert()

# Header cell for the newly inserted column.
$ws.Range("V1").Value = "Si autre situation pro"

# New column adopts the neighbouring (former column U) header/data look,
# minus the border - drop the border explicitly so the engine records a
# distinct (border-less) style for these cells, same as column U's style
# but without the grid border.
$ws.Range("V1").Borders.LineStyle = -4142
$ws.Range("V2:V4").Borders.LineStyle = -4142

# The hidden _FilterDatabase defined name covered column A:AY on row 1;
# extend it by one column to A:AZ to keep covering the full header row.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Feuil1!`$A`$1:`$AZ`$1"
    }
}

# Move the active selection onto the newly inserted column's first data
# cell, mirroring the author's final cursor position.
$ws.Range("V2").Select() | Out-Null
